# adding averages and more checks
# Update the "Training Dashboard" sheet: PERIOD TO EXPIRE (H) shrinks by 8
# days and LAST UPDATE (I) moves from 08-Sep-2025 to 16-Sep-2025 for every
# training row.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Training Dashboard")

$ws1.Range("H3").Value  = 334
$ws1.Range("H4").Value  = 638
$ws1.Range("H5").Value  = 364
$ws1.Range("H6").Value  = 435
$ws1.Range("H7").Value  = 423
$ws1.Range("H8").Value  = 86
$ws1.Range("H9").Value  = -343
$ws1.Range("H10").Value = 182

# Column I holds dates stored as plain text ("dd-mmm-yyyy"). Force the
# number format to Text first so Excel doesn't silently reinterpret the
# string as a real date serial when we write it back.
$ws1.Range("I3:I10").NumberFormat = "@"
$ws1.Range("I3").Value  = "16-Sep-2025"
$ws1.Range("I4").Value  = "16-Sep-2025"
$ws1.Range("I5").Value  = "16-Sep-2025"
$ws1.Range("I6").Value  = "16-Sep-2025"
$ws1.Range("I7").Value  = "16-Sep-2025"
$ws1.Range("I8").Value  = "16-Sep-2025"
$ws1.Range("I9").Value  = "16-Sep-2025"
$ws1.Range("I10").Value = "16-Sep-2025"

# Update the "Exam Dashboard" sheet: widen the COMMENTS column and replace
# the generic "OK" comment with a more descriptive check result.
$ws2 = $wb.Worksheets.Item("Exam Dashboard")
$ws2.Range("E1").ColumnWidth = 14.1666666667

$ws2.Range("E3").Value  = "date is valid"
$ws2.Range("E4").Value  = "date is valid"
$ws2.Range("E5").Value  = "date is valid"
$ws2.Range("E6").Value  = "date is valid"
$ws2.Range("E7").Value  = "date is valid"
$ws2.Range("E8").Value  = "date is valid"
$ws2.Range("E9").Value  = "date is valid"
$ws2.Range("E10").Value = "date is valid"
